$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D4").Value = -7.695800000000004
$ws.Range("D6").Value = -7.664
$ws.Range("D7").Value = -7.690600000000002
$ws.Range("D8").Value = -8.298600000000002
$ws.Range("D16").Value = -8.048599999999995
$ws.Range("D20").Value = -8.576899999999995
$ws.Range("D21").Value = -8.152799999999996
$ws.Range("D28").Value = -8.370699999999998
$ws.Range("D29").Value = -7.247300000000005
$ws.Range("D30").Value = -6.999199999999995
$ws.Range("D32").Value = -7.504099999999994
$ws.Range("D40").Value = -8.735699999999994
$ws.Range("D46").Value = -7.479999999999998
$ws.Range("D51").Value = -8.1488
$ws.Range("D52").Value = -7.754000000000002
$ws.Range("D57").Value = -8.075699999999999
$ws.Range("D59").Value = -8.239299999999998
$ws.Range("D62").Value = -8.509099999999995
$ws.Range("D66").Value = -6.930799999999998
$ws.Range("D73").Value = -7.736299999999989
$ws.Range("D74").Value = -8.289300000000008
$ws.Range("D77").Value = -6.136000000000001
$ws.Range("D92").Value = -6.245700000000006
$ws.Range("D100").Value = -7.795100000000002
